$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-9 with corrected/simulated values (columns B:F),
# per bugfixed evaluation and simulated rt_data.

# Row 2 (Q0)
$ws.Range("B2").Value = -0.3950010431191253
$ws.Range("C2").Value = 2.046336902374465
$ws.Range("D2").Value = 5.798338224185904
$ws.Range("E2").Value = 2.407973883617907
$ws.Range("F2").Value = 2.46502260251285

# Row 3 (Q1)
$ws.Range("B3").Value = -0.05523834182961166
$ws.Range("C3").Value = 1.856073131361144
$ws.Range("D3").Value = 4.905587478437071
$ws.Range("E3").Value = 2.214856085265377
$ws.Range("F3").Value = 2.304578245080446

# Row 4 (Q2)
$ws.Range("B4").Value = -0.290471866325093
$ws.Range("C4").Value = 1.739315905158552
$ws.Range("D4").Value = 4.447440393765628
$ws.Range("E4").Value = 2.108895538846253
$ws.Range("F4").Value = 2.181675708417124

# Row 5 (Q3)
$ws.Range("B5").Value = -0.374175682380495
$ws.Range("C5").Value = 1.675777845656498
$ws.Range("D5").Value = 3.969686682774517
$ws.Range("E5").Value = 1.992407258261854
$ws.Range("F5").Value = 2.052473426292913

# Row 6 (Q4)
$ws.Range("B6").Value = -0.589146841836848
$ws.Range("C6").Value = 1.539691385723428
$ws.Range("D6").Value = 3.478629428342508
$ws.Range("E6").Value = 1.865108422677488
$ws.Range("F6").Value = 1.865337451477517

# Row 7 (Q5)
$ws.Range("B7").Value = -0.1732891586448268
$ws.Range("C7").Value = 1.168271714243534
$ws.Range("D7").Value = 2.179932902258912
$ws.Range("E7").Value = 1.476459583686229
$ws.Range("F7").Value = 1.555198296351454
$ws.Range("G7").Value = 9

# Row 8 (Q6)
$ws.Range("B8").Value = -0.2018953961280825
$ws.Range("C8").Value = 1.273161211529873
$ws.Range("D8").Value = 2.752380755634883
$ws.Range("E8").Value = 1.659030064716997
$ws.Range("F8").Value = 1.803868843787874
$ws.Range("G8").Value = 6

# Row 9 (Q7)
$ws.Range("B9").Value = 1.036650915497801
$ws.Range("C9").Value = 1.036650915497801
$ws.Range("D9").Value = 1.534553857341397
$ws.Range("E9").Value = 1.238771107727895
$ws.Range("F9").Value = 0.8305799811628327
$ws.Range("G9").Value = 3

# New row 10 (Q8)
$ws.Range("A10").Value = "Q8"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").Value = -0.1948090203170254
$ws.Range("C10").Value = 0.1948090203170254
$ws.Range("D10").Value = 0.03795055439687921
$ws.Range("E10").Value = 0.1948090203170254
$ws.Range("G10").Value = 1
